{"js": "const body = context.document.body;\n\n{\n  const results = body.search(\"2025-04-28 Monday\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  results.items[0].insertText(\"2025-04-29 Tuesday\", Word.InsertLocation.replace);\n}\n{\n  const results = body.search(\"77\u00d769=\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  results.items[0].insertText(\"29\u00d725=\", Word.InsertLocation.replace);\n}\n{\n  const results = body.search(\"32\u00d711=\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  results.items[0].insertText(\"74\u00d767=\", Word.InsertLocation.replace);\n}\n{\n  const results = body.search(\"51\u00d797=\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  results.items[0].insertText(\"11\u00d771=\", Word.InsertLocation.replace);\n}\n{\n  const results = body.search(\"68\u00d786=\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  results.items[0].insertText(\"56\u00d712=\", Word.InsertLocation.replace);\n}\n{\n  const results = body.search(\"17\u00d738=\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  results.items[0].insertText(\"70\u00d721=\", Word.InsertLocation.replace);\n}\n{\n  const results = body.search(\"23\u00d757=\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  results.items[0].insertText(\"54\u00d712=\", Word.InsertLocation.replace);\n}\n{\n  const results = body.search(\"47\u00d774=\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  results.items[0].insertText(\"69\u00d725=\", Word.InsertLocation.replace);\n}\n{\n  const results = body.search(\"51\u00d747=\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  results.items[0].insertText(\"19\u00d755=\", Word.InsertLocation.replace);\n}\n{\n  const results = body.search(\"15\u00d739=\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  results.items[0].insertText(\"92\u00d770=\", Word.InsertLocation.replace);\n}\n{\n  const results = body.search(\"34\u00d733=\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  results.items[0].insertText(\"43\u00d772=\", Word.InsertLocation.replace);\n}\n{\n  const results = body.search(\"43\u00d732=\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  results.items[0].insertText(\"55\u00d747=\", Word.InsertLocation.replace);\n}\n{\n  const results = body.search(\"48\u00d714=\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  results.items[0].insertText(\"40\u00d799=\", Word.InsertLocation.replace);\n}\n{\n  const results = body.search(\"49\u00d747=\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  results.items[0].insertText(\"53\u00d732=\", Word.InsertLocation.replace);\n}\n{\n  const results = body.search(\"46\u00d764=\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  results.items[0].insertText(\"92\u00d784=\", Word.InsertLocation.replace);\n}\n{\n  const results = body.search(\"48\u00d757=\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  results.items[0].insertText(\"97\u00d726=\", Word.InsertLocation.replace);\n}\n{\n  const results = body.search(\"71\u00d714=\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  results.items[0].insertText(\"83\u00d798=\", Word.InsertLocation.replace);\n}\n{\n  const results = body.search(\"67\u00d712=\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  results.items[0].insertText(\"55\u00d757=\", Word.InsertLocation.replace);\n}\n{\n  const results = body.search(\"82\u00d715=\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  results.items[0].insertText(\"57\u00d764=\", Word.InsertLocation.replace);\n}\n{\n  const results = body.search(\"92\u00d757=\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  results.items[0].insertText(\"83\u00d742=\", Word.InsertLocation.replace);\n}\n{\n  const results = body.search(\"89\u00d755=\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  results.items[0].insertText(\"45\u00d731=\", Word.InsertLocation.replace);\n}\n{\n  const results = body.search(\"13\u00d754=\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  results.items[0].insertText(\"49\u00d725=\", Word.InsertLocation.replace);\n}\n{\n  const results = body.search(\"42\u00d711=\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  results.items[0].insertText(\"21\u00d767=\", Word.InsertLocation.replace);\n}\n{\n  const results = body.search(\"32\u00d789=\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  results.items[0].insertText(\"49\u00d793=\", Word.InsertLocation.replace);\n}\n{\n  const results = body.search(\"17\u00d735=\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  results.items[0].insertText(\"22\u00d730=\", Word.InsertLocation.replace);\n}\n{\n  const results = body.search(\"61\u00d779=\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  results.items[0].insertText(\"73\u00d775=\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.Execute(\"2025-04-28 Monday\", $false, $false, $false, $false, $false, $true, 1, $false, \"2025-04-29 Tuesday\", 2) | Out-Null\n$find = $d.Content.Find\n$find.Execute(\"77\u00d769=\", $false, $false, $false, $false, $false, $true, 1, $false, \"29\u00d725=\", 2) | Out-Null\n$find = $d.Content.Find\n$find.Execute(\"32\u00d711=\", $false, $false, $false, $false, $false, $true, 1, $false, \"74\u00d767=\", 2) | Out-Null\n$find = $d.Content.Find\n$find.Execute(\"51\u00d797=\", $false, $false, $false, $false, $false, $true, 1, $false, \"11\u00d771=\", 2) | Out-Null\n$find = $d.Content.Find\n$find.Execute(\"68\u00d786=\", $false, $false, $false, $false, $false, $true, 1, $false, \"56\u00d712=\", 2) | Out-Null\n$find = $d.Content.Find\n$find.Execute(\"17\u00d738=\", $false, $false, $false, $false, $false, $true, 1, $false, \"70\u00d721=\", 2) | Out-Null\n$find = $d.Content.Find\n$find.Execute(\"23\u00d757=\", $false, $false, $false, $false, $false, $true, 1, $false, \"54\u00d712=\", 2) | Out-Null\n$find = $d.Content.Find\n$find.Execute(\"47\u00d774=\", $false, $false, $false, $false, $false, $true, 1, $false, \"69\u00d725=\", 2) | Out-Null\n$find = $d.Content.Find\n$find.Execute(\"51\u00d747=\", $false, $false, $false, $false, $false, $true, 1, $false, \"19\u00d755=\", 2) | Out-Null\n$find = $d.Content.Find\n$find.Execute(\"15\u00d739=\", $false, $false, $false, $false, $false, $true, 1, $false, \"92\u00d770=\", 2) | Out-Null\n$find = $d.Content.Find\n$find.Execute(\"34\u00d733=\", $false, $false, $false, $false, $false, $true, 1, $false, \"43\u00d772=\", 2) | Out-Null\n$find = $d.Content.Find\n$find.Execute(\"43\u00d732=\", $false, $false, $false, $false, $false, $true, 1, $false, \"55\u00d747=\", 2) | Out-Null\n$find = $d.Content.Find\n$find.Execute(\"48\u00d714=\", $false, $false, $false, $false, $false, $true, 1, $false, \"40\u00d799=\", 2) | Out-Null\n$find = $d.Content.Find\n$find.Execute(\"49\u00d747=\", $false, $false, $false, $false, $false, $true, 1, $false, \"53\u00d732=\", 2) | Out-Null\n$find = $d.Content.Find\n$find.Execute(\"46\u00d764=\", $false, $false, $false, $false, $false, $true, 1, $false, \"92\u00d784=\", 2) | Out-Null\n$find = $d.Content.Find\n$find.Execute(\"48\u00d757=\", $false, $false, $false, $false, $false, $true, 1, $false, \"97\u00d726=\", 2) | Out-Null\n$find = $d.Content.Find\n$find.Execute(\"71\u00d714=\", $false, $false, $false, $false, $false, $true, 1, $false, \"83\u00d798=\", 2) | Out-Null\n$find = $d.Content.Find\n$find.Execute(\"67\u00d712=\", $false, $false, $false, $false, $false, $true, 1, $false, \"55\u00d757=\", 2) | Out-Null\n$find = $d.Content.Find\n$find.Execute(\"82\u00d715=\", $false, $false, $false, $false, $false, $true, 1, $false, \"57\u00d764=\", 2) | Out-Null\n$find = $d.Content.Find\n$find.Execute(\"92\u00d757=\", $false, $false, $false, $false, $false, $true, 1, $false, \"83\u00d742=\", 2) | Out-Null\n$find = $d.Content.Find\n$find.Execute(\"89\u00d755=\", $false, $false, $false, $false, $false, $true, 1, $false, \"45\u00d731=\", 2) | Out-Null\n$find = $d.Content.Find\n$find.Execute(\"13\u00d754=\", $false, $false, $false, $false, $false, $true, 1, $false, \"49\u00d725=\", 2) | Out-Null\n$find = $d.Content.Find\n$find.Execute(\"42\u00d711=\", $false, $false, $false, $false, $false, $true, 1, $false, \"21\u00d767=\", 2) | Out-Null\n$find = $d.Content.Find\n$find.Execute(\"32\u00d789=\", $false, $false, $false, $false, $false, $true, 1, $false, \"49\u00d793=\", 2) | Out-Null\n$find = $d.Content.Find\n$find.Execute(\"17\u00d735=\", $false, $false, $false, $false, $false, $true, 1, $false, \"22\u00d730=\", 2) | Out-Null\n$find = $d.Content.Find\n$find.Execute(\"61\u00d779=\", $false, $false, $false, $false, $false, $true, 1, $false, \"73\u00d775=\", 2) | Out-Null\n"}
